$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exhibitor")
$ws.Name = "expo"
